$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.596.41'
$ws.Range("E2").Value = '  +0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.478.36'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9602'
$ws.Range("E5").Value = '  +5.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '280.73'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3665'
$ws.Range("E7").Value = '  -1.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3071'
$ws.Range("E8").Value = '  -3.79%  '

$ws.Range("E9").Value = '  -0.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.063'
$ws.Range("E10").Value = '  +0.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06680'
$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.530'
$ws.Range("E13").Value = '  -0.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.09'
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.224'
$ws.Range("E15").Value = '  -0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9610'
$ws.Range("E16").Value = '  +4.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001035'
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.478.30'
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05959'
$ws.Range("E19").Value = '  +3.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.06'
$ws.Range("E20").Value = '  -2.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.510'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.46'
$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  -1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.267'
$ws.Range("E24").Value = '  -1.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.622.72'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.51'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.113'
$ws.Range("E27").Value = '  -8.42%  '

$ws.Range("E28").Value = '  -1.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.637.90'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.86'
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.970'
$ws.Range("E31").Value = '  +0.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.034'
$ws.Range("E32").Value = '  -4.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8159'
$ws.Range("E33").Value = '  -4.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07976'
$ws.Range("E34").Value = '  +2.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.533'
$ws.Range("E35").Value = '  -0.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.220'
$ws.Range("E36").Value = '  +3.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05806'
$ws.Range("E37").Value = '  -5.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.741'
$ws.Range("E38").Value = '  -3.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02051'
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9608'
$ws.Range("E40").Value = '  +2.55%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.40'
$ws.Range("E41").Value = '  -2.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1879'
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.490'
$ws.Range("E43").Value = '  +1.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5319'
$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.544'
$ws.Range("E45").Value = '  -1.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.31'
$ws.Range("E46").Value = '  -1.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.35'
$ws.Range("E47").Value = '  -4.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5206'
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.824'
$ws.Range("E49").Value = '  -0.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06494'
$ws.Range("E50").Value = '  +0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9916'
$ws.Range("E51").Value = '  -0.14%  '
